$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.533.80'
$ws.Range("E2").Value = '  +3.33%  '
$ws.Range("D3").Value = '3.071.41'
$ws.Range("E3").Value = '  +2.41%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '551.86'
$ws.Range("E5").Value = '  +2.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.68'
$ws.Range("E6").Value = '  +5.85%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '3.068.29'
$ws.Range("E8").Value = '  +2.39%  '
$ws.Range("E9").Value = '  +1.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.52'
$ws.Range("E10").Value = '  +6.24%  '
$ws.Range("E11").Value = '  +2.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.459'
$ws.Range("E12").Value = '  +2.36%  '
$ws.Range("E13").Value = '  +2.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.00'
$ws.Range("E14").Value = '  +2.92%  '
$ws.Range("D15").Value = '3.564.47'
$ws.Range("E15").Value = '  +2.34%  '
$ws.Range("D16").Value = '63.497.53'
$ws.Range("D17").Value = '3.070.97'
$ws.Range("E17").Value = '  +2.39%  '
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("E19").Value = '  +2.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '485.61'
$ws.Range("E20").Value = '  +3.86%  '
$ws.Range("E21").Value = '  +5.23%  '
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("E23").Value = '  +5.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.21'
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.78'
$ws.Range("E25").Value = '  +6.43%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.78'
$ws.Range("E27").Value = '  +3.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.94'
$ws.Range("E28").Value = '  +2.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.03'
$ws.Range("E29").Value = '  +7.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.29'
$ws.Range("E31").Value = '  +2.57%  '
$ws.Range("E32").Value = '  +1.58%  '
$ws.Range("E33").Value = '  +7.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.71'
$ws.Range("E34").Value = '  +3.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.62'
$ws.Range("E35").Value = '  +1.51%  '
$ws.Range("E36").Value = '  +1.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '467.35'
$ws.Range("E38").Value = '  +4.95%  '
$ws.Range("E39").Value = '  +3.88%  '
$ws.Range("D40").Value = '3.047.66'
$ws.Range("E40").Value = '  -3.71%  '
$ws.Range("E41").Value = '  -1.12%  '
$ws.Range("E42").Value = '  +1.44%  '
$ws.Range("E43").Value = '  +5.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.81'
$ws.Range("E44").Value = '  +2.76%  '
$ws.Range("E45").Value = '  +4.96%  '
$ws.Range("E47").Value = '  +3.16%  '
$ws.Range("E48").Value = '  +2.55%  '
$ws.Range("E49").Value = '  +3.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '116.99'
$ws.Range("E50").Value = '  -1.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.10'
$ws.Range("E51").Value = '  +4.39%  '
